# TC04_C3DC_phs003164_AnatomicSite-C421BoneMarrow.xlsx
# "Updated C3DC Regression and Smoke suites"
#
# The TreatmentTab query stored in Sheet1!B5 wrapped its REPLACE(...) call in
# a redundant CONCAT(...) - fix the "Treatment Agent" column expression so it
# simply is REPLACE(trt.treatment_agent, ';', ', ').

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentCell = $ws.Range("B5")
$oldQuery = [string]$treatmentCell.Value2

$badExpr  = "CONCAT(REPLACE(trt.treatment_agent, ';', ', '))"
$goodExpr = "REPLACE(trt.treatment_agent, ';', ', ')"
$newQuery = $oldQuery.Replace($badExpr, $goodExpr)

$treatmentCell.Value = $newQuery

# Restore the plain B2 selection (no scrolled topLeftCell) that the saved
# workbook ends up with.
[void]$ws.Range("B2").Select()
